$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 6676
$ws.Range('K3').Value = 6894
$ws.Range('D4').Value = 1978
$ws.Range('K4').Value = 1431
$ws.Range('K5').Value = 498
$ws.Range('K6').Value = 7570
$ws.Range('D7').Value = 28168
$ws.Range('K7').Value = 23069

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K3').Value = 79
$ws.Range('K4').Value = 14
$ws.Range('K7').Value = 293

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 417
$ws.Range('K3').Value = 461
$ws.Range('K4').Value = 87
$ws.Range('K6').Value = 500
$ws.Range('K7').Value = 1511

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K2').Value = 174
$ws.Range('K3').Value = 177
$ws.Range('K7').Value = 499

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 354
$ws.Range('K6').Value = 315
$ws.Range('K7').Value = 999

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 224
$ws.Range('K3').Value = 258
$ws.Range('K7').Value = 781

$ws = $wb.Worksheets.Item('New City')
$ws.Range('K3').Value = 132
$ws.Range('K6').Value = 196
$ws.Range('K7').Value = 539

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K2').Value = 102
$ws.Range('K7').Value = 391

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K2').Value = 201
$ws.Range('K7').Value = 698
$ws.Range('K8').Value = 1511
$ws.Range('K10').Value = 134
$ws.Range('K11').Value = 425
$ws.Range('K12').Value = 41
$ws.Range('K18').Value = 154
$ws.Range('K19').Value = 674
$ws.Range('K20').Value = 557
$ws.Range('K27').Value = 216
$ws.Range('K29').Value = 1251
$ws.Range('K31').Value = 254
$ws.Range('K33').Value = 999
$ws.Range('K34').Value = 130
$ws.Range('K36').Value = 293
$ws.Range('K37').Value = 781
$ws.Range('K42').Value = 851
$ws.Range('K44').Value = 192
$ws.Range('K48').Value = 293
$ws.Range('K51').Value = 289
$ws.Range('K52').Value = 612
$ws.Range('K53').Value = 293
$ws.Range('K55').Value = 248
$ws.Range('D63').Value = 357
$ws.Range('K63').Value = 59
$ws.Range('K65').Value = 539
$ws.Range('K67').Value = 901
$ws.Range('K68').Value = 62
$ws.Range('K71').Value = 70
$ws.Range('K76').Value = 310
$ws.Range('K77').Value = 157
$ws.Range('K78').Value = 262
$ws.Range('K79').Value = 577
$ws.Range('K80').Value = 84
$ws.Range('K83').Value = 499
$ws.Range('K84').Value = 186
$ws.Range('K85').Value = 1065
$ws.Range('K86').Value = 141
$ws.Range('K88').Value = 247
$ws.Range('K89').Value = 346
$ws.Range('K90').Value = 219
$ws.Range('K91').Value = 274
$ws.Range('K92').Value = 85
$ws.Range('K94').Value = 309
$ws.Range('K96').Value = 244
$ws.Range('K98').Value = 116
$ws.Range('K99').Value = 391
$ws.Range('D101').Value = 28168
$ws.Range('K101').Value = 23069

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K6').Value = 89
$ws.Range('K7').Value = 254

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('K2').Value = 245
$ws.Range('K3').Value = 329
$ws.Range('K6').Value = 257
$ws.Range('K7').Value = 901

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 76
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 186

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 353
$ws.Range('K3').Value = 443
$ws.Range('K4').Value = 60
$ws.Range('K5').Value = 31
$ws.Range('K7').Value = 1251

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K5').Value = 2
$ws.Range('K7').Value = 293

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 198
$ws.Range('K3').Value = 202
$ws.Range('K4').Value = 32
$ws.Range('K7').Value = 674

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('K6').Value = 76
$ws.Range('K7').Value = 192

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K2').Value = 72
$ws.Range('K4').Value = 20
$ws.Range('K7').Value = 310

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 229
$ws.Range('K6').Value = 315
$ws.Range('K7').Value = 851

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('K2').Value = 43
$ws.Range('K7').Value = 134

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('K2').Value = 78
$ws.Range('K7').Value = 262

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K2').Value = 76
$ws.Range('K3').Value = 72
$ws.Range('K7').Value = 248

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K3').Value = 49
$ws.Range('K7').Value = 244

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K3').Value = 130
$ws.Range('K7').Value = 274

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K3').Value = 186
$ws.Range('K4').Value = 36
$ws.Range('K6').Value = 143
$ws.Range('K7').Value = 577

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K2').Value = 194
$ws.Range('K4').Value = 27
$ws.Range('K7').Value = 557

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('K3').Value = 51
$ws.Range('K7').Value = 154

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K3').Value = 89
$ws.Range('K7').Value = 293

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 228
$ws.Range('K3').Value = 226
$ws.Range('K7').Value = 698

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('K3').Value = 35
$ws.Range('K7').Value = 130

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('K2').Value = 78
$ws.Range('K4').Value = 24
$ws.Range('K6').Value = 140
$ws.Range('K7').Value = 309

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('K6').Value = 68
$ws.Range('K7').Value = 116

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K2').Value = 147
$ws.Range('K7').Value = 425

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('K6').Value = 64
$ws.Range('K7').Value = 201

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('K3').Value = 17
$ws.Range('K7').Value = 85

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K2').Value = 64
$ws.Range('K6').Value = 99
$ws.Range('K7').Value = 247

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K3').Value = 108
$ws.Range('K7').Value = 346

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range('K2').Value = 56
$ws.Range('K7').Value = 216

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('K2').Value = 24
$ws.Range('K7').Value = 141

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range('K2').Value = 82
$ws.Range('K7').Value = 219

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 77
$ws.Range('K7').Value = 289

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('K3').Value = 13
$ws.Range('K7').Value = 62

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 351
$ws.Range('K3').Value = 371
$ws.Range('K7').Value = 1065

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range('K3').Value = 25
$ws.Range('K6').Value = 19
$ws.Range('K7').Value = 70

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('K2').Value = 65
$ws.Range('K7').Value = 157

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range('K6').Value = 39
$ws.Range('K7').Value = 84

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('K2').Value = 164
$ws.Range('K4').Value = 34
$ws.Range('K6').Value = 222
$ws.Range('K7').Value = 612

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('K6').Value = 15
$ws.Range('K7').Value = 41
